$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1222.7059
$ws.Range("I80").Value = 1886.5714
$ws.Range("J80").Value = 758
$ws.Range("K80").Value = 5659.7142
$ws.Range("L80").Value = 2274
$ws.Range("M80").Value = -4661.7142
$ws.Range("N80").Value = -4270

$ws.Range("H83").Value = 1222.7059
$ws.Range("I83").Value = 1886.5714
$ws.Range("J83").Value = 758
$ws.Range("K83").Value = 16979.1426
$ws.Range("L83").Value = 6822
$ws.Range("M83").Value = -11987.1426
$ws.Range("N83").Value = -16806

$ws.Range("H111").Value = 925.8
$ws.Range("I111").Value = 825.6
$ws.Range("J111").Value = 1026
$ws.Range("K111").Value = 2476.8
$ws.Range("L111").Value = 3078
$ws.Range("M111").Value = 590.1999999999998
$ws.Range("N111").Value = -9212

$ws.Range("H113").Value = 2142.8572
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254

$ws.Range("H125").Value = 1202.125
$ws.Range("I125").Value = 1690.25
$ws.Range("J125").Value = 714
$ws.Range("K125").Value = 15212.25
$ws.Range("L125").Value = 6426
$ws.Range("M125").Value = -12752.25
$ws.Range("N125").Value = -11346

$ws.Range("H129").Value = 1229.8077
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1229.8077
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3689.4231
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13689.4231

$ws.Range("H131").Value = 4207.9287
$ws.Range("I131").Value = 1255.091
$ws.Range("J131").Value = 15035
$ws.Range("K131").Value = 3765.273
$ws.Range("L131").Value = 45105
$ws.Range("M131").Value = 1274.727
$ws.Range("N131").Value = -55185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10566.243
$ws.Range("I32").Value = 9159.709999999999
$ws.Range("J32").Value = 17833.334
$ws.Range("K32").Value = 9159.709999999999
$ws.Range("L32").Value = 17833.334
$ws.Range("M32").Value = -8872.709999999999
$ws.Range("N32").Value = -18407.334

$ws.Range("H61").Value = 20835366
$ws.Range("I61").Value = 26317816
$ws.Range("J61").Value = 2060.8
$ws.Range("K61").Value = 26317816
$ws.Range("L61").Value = 2060.8
$ws.Range("M61").Value = -26317604
$ws.Range("N61").Value = -2484.8

$ws.Range("H110").Value = 1324.4
$ws.Range("I110").Value = 1186.875
$ws.Range("J110").Value = 1874.5
$ws.Range("K110").Value = 1186.875
$ws.Range("L110").Value = 1874.5
$ws.Range("M110").Value = 858.125
$ws.Range("N110").Value = -5964.5

$ws.Range("H133").Value = 51931.715
$ws.Range("J133").Value = 51931.715
$ws.Range("L133").Value = 51931.715
$ws.Range("N133").Value = -56991.715

$ws.Range("H136").Value = 20835366
$ws.Range("I136").Value = 26317816
$ws.Range("J136").Value = 2060.8
$ws.Range("K136").Value = 78953448
$ws.Range("L136").Value = 6182.400000000001
$ws.Range("M136").Value = -78950898
$ws.Range("N136").Value = -11282.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1369.3572
$ws.Range("I107").Value = 1365.7778
$ws.Range("J107").Value = 1375.8
$ws.Range("K107").Value = 1365.7778
$ws.Range("L107").Value = 1375.8
$ws.Range("M107").Value = 554.2221999999999
$ws.Range("N107").Value = -5215.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1200.1666
$ws.Range("I16").Value = 967
$ws.Range("J16").Value = 1433.3334
$ws.Range("K16").Value = 967
$ws.Range("L16").Value = 1433.3334
$ws.Range("M16").Value = -680
$ws.Range("N16").Value = -2007.3334

$ws.Range("H31").Value = 7096355.5
$ws.Range("I31").Value = 4809.0884
$ws.Range("J31").Value = 25643478
$ws.Range("K31").Value = 4809.0884
$ws.Range("L31").Value = 25643478
$ws.Range("M31").Value = -4514.0884
$ws.Range("N31").Value = -25644068

$ws.Range("H34").Value = 7096355.5
$ws.Range("I34").Value = 4809.0884
$ws.Range("J34").Value = 25643478
$ws.Range("K34").Value = 4809.0884
$ws.Range("L34").Value = 25643478
$ws.Range("M34").Value = -4607.0884
$ws.Range("N34").Value = -25643882

$ws.Range("H107").Value = 690.7778
$ws.Range("I107").Value = 764.625
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 764.625
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1155.375
$ws.Range("N107").Value = -3940

$ws.Range("H113").Value = 1200.1666
$ws.Range("I113").Value = 967
$ws.Range("J113").Value = 1433.3334
$ws.Range("K113").Value = 967
$ws.Range("L113").Value = 1433.3334
$ws.Range("M113").Value = 1203
$ws.Range("N113").Value = -5773.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1262.25
$ws.Range("I69").Value = 724.7143
$ws.Range("J69").Value = 2516.5
$ws.Range("K69").Value = 2174.1429
$ws.Range("L69").Value = 7549.5
$ws.Range("M69").Value = -1363.1429
$ws.Range("N69").Value = -9171.5

$ws.Range("H72").Value = 1262.25
$ws.Range("I72").Value = 724.7143
$ws.Range("J72").Value = 2516.5
$ws.Range("K72").Value = 6522.428699999999
$ws.Range("L72").Value = 22648.5
$ws.Range("M72").Value = -2466.428699999999
$ws.Range("N72").Value = -30760.5

$ws.Range("H113").Value = 500.1837
$ws.Range("I113").Value = 437.02777
$ws.Range("J113").Value = 675.0769
$ws.Range("K113").Value = 1311.08331
$ws.Range("L113").Value = 2025.2307
$ws.Range("M113").Value = 858.91669
$ws.Range("N113").Value = -6365.2307

$ws.Range("H132").Value = 1716.5
$ws.Range("I132").Value = 833
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 7497
$ws.Range("L132").Value = 23400
$ws.Range("M132").Value = -4967
$ws.Range("N132").Value = -28460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8204079
$ws.Range("I80").Value = 12823011
$ws.Range("J80").Value = 1883435.1
$ws.Range("K80").Value = 12823011
$ws.Range("L80").Value = 1883435.1
$ws.Range("M80").Value = -12822013
$ws.Range("N80").Value = -1885431.1

$ws.Range("H83").Value = 8204079
$ws.Range("I83").Value = 12823011
$ws.Range("J83").Value = 1883435.1
$ws.Range("K83").Value = 64115055
$ws.Range("L83").Value = 9417175.5
$ws.Range("M83").Value = -64110063
$ws.Range("N83").Value = -9427159.5

$ws.Range("H102").Value = 3622.2068
$ws.Range("I102").Value = 3817.1538
$ws.Range("J102").Value = 1932.6666
$ws.Range("K102").Value = 3817.1538
$ws.Range("L102").Value = 1932.6666
$ws.Range("M102").Value = -2195.1538
$ws.Range("N102").Value = -5176.6666

$ws.Range("H122").Value = 3176586
$ws.Range("I122").Value = 4445918.5
$ws.Range("K122").Value = 13337755.5
$ws.Range("M122").Value = -13335305.5

$ws.Range("H132").Value = 2982.5642
$ws.Range("I132").Value = 2517.3333
$ws.Range("J132").Value = 4533.3335
$ws.Range("K132").Value = 7551.999899999999
$ws.Range("L132").Value = 13600.0005
$ws.Range("M132").Value = -5021.999899999999
$ws.Range("N132").Value = -18660.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1555.3334
$ws.Range("I68").Value = 1960.6154
$ws.Range("J68").Value = 1076.3636
$ws.Range("K68").Value = 1960.6154
$ws.Range("L68").Value = 1076.3636
$ws.Range("M68").Value = -1211.6154
$ws.Range("N68").Value = -2574.3636

$ws.Range("H71").Value = 1555.3334
$ws.Range("I71").Value = 1960.6154
$ws.Range("J71").Value = 1076.3636
$ws.Range("K71").Value = 9803.076999999999
$ws.Range("L71").Value = 5381.817999999999
$ws.Range("M71").Value = -6059.076999999999
$ws.Range("N71").Value = -12869.818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1245.5385
$ws.Range("I107").Value = 1484.6
$ws.Range("J107").Value = 448.66666
$ws.Range("K107").Value = 4453.799999999999
$ws.Range("L107").Value = 1345.99998
$ws.Range("M107").Value = -2533.799999999999
$ws.Range("N107").Value = -5185.999980000001

$ws.Range("H122").Value = 2103.6191
$ws.Range("I122").Value = 2104
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 6312
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -3862
$ws.Range("N122").Value = -11200

$ws.Range("H126").Value = 1986.4286
$ws.Range("I126").Value = 1309.3478
$ws.Range("J126").Value = 5101
$ws.Range("K126").Value = 3928.0434
$ws.Range("L126").Value = 15303
$ws.Range("M126").Value = -1458.0434
$ws.Range("N126").Value = -20243
